$wb = $excel.ActiveWorkbook

# --- Properties sheet: add 3 new GeneratingUnit properties (cable, operator, creationTime) ---
$props = $wb.Worksheets.Item("Properties")

# Copy the formatting (fill + border style) of the last existing data row (row 6)
# down into the three new rows (7, 8, 9) before filling in values.
$props.Range("A6:P6").Copy()
$props.Range("A7:P9").PasteSpecial(-4122)

$newProps = @("cable", "operator", "creationTime")
$row = 7
foreach ($propName in $newProps) {
    $props.Range("A" + $row).Value = "GeneratingUnit"
    $props.Range("B" + $row).Value = $propName
    $props.Range("F" + $row).Value = "text"
    $props.Range("G" + $row).Value = $false
    $props.Range("H" + $row).Value = $false
    $props.Range("K" + $row).Value = "GeneratingUnit"
    $props.Range("L" + $row).Value = $propName
    $props.Range("M" + $row).Value = "GeneratingUnit"
    $props.Range("N" + $row).Value = $propName
    $row = $row + 1
}

# --- Switch the active sheet/selection from Containers to Properties ---
$props.Activate()
$props.Range("L10").Select()
